# "Mas ejercicios y correciones." -- add a 4th exercise/attendance column (D)
# to the existing 3 (B, C) and backfill it for every student row, plus a
# numeric header cell matching the existing B1/C1 header pattern.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.ActiveSheet
$win = $excel.ActiveWindow

# New column D should look just like column C (same centered style), so
# clone C1:C20's formatting onto D1:D20 before filling in values.
$ws.Range("C1:C20").Copy()
[void]$ws.Range("D1:D20").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Header cell for the new exercise column.
$ws.Range("D1").Value = 3

# Attendance marks ("x") for the new column -- mirrors the rows that already
# carry a mark in column C, minus row 10, plus the new marks on rows 18/20
# (row 17/19 stay blank in column D).
$ws.Range("D2").Value  = "x"
$ws.Range("D3").Value  = "x"
$ws.Range("D4").Value  = "x"
$ws.Range("D5").Value  = "x"
$ws.Range("D6").Value  = "x"
$ws.Range("D9").Value  = "x"
$ws.Range("D11").Value = "x"
$ws.Range("D18").Value = "x"
$ws.Range("D20").Value = "x"

# Rows 7, 8, 10, 12-17, 19 keep an empty (but styled) D cell -- already
# created by the PasteSpecial above, nothing further to set.

# View state: zoomed in and scrolled/selected near the bottom of the sheet,
# with the selection now anchored on the newly-extended column.
$win.Zoom = 162
[void]$ws.Range("D21").Select()
